# "add ssd300 voc and coco config"
# Expand the experiments table from 7 columns x 3 rows (A1:G3) to
# 13 columns x 6 rows (A1:M6), restyle header/body fonts & alignment,
# resize the new columns, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths (A:M)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 13.6640625
$ws.Columns.Item(2).ColumnWidth  = 17.33203125
$ws.Columns.Item(3).ColumnWidth  = 13.6640625
$ws.Columns.Item(4).ColumnWidth  = 16.109375
$ws.Columns.Item(5).ColumnWidth  = 13.21875
$ws.Columns.Item(6).ColumnWidth  = 13.109375
$ws.Columns.Item(7).ColumnWidth  = 8.33203125
$ws.Columns.Item(8).ColumnWidth  = 18.5546875
$ws.Columns.Item(9).ColumnWidth  = 14.44140625
$ws.Columns.Item(10).ColumnWidth = 17.33203125
$ws.Columns.Item(11).ColumnWidth = 17.33203125
$ws.Columns.Item(12).ColumnWidth = 17.33203125
$ws.Columns.Item(13).ColumnWidth = 13.21875

# ---------------------------------------------------------------------------
# 2. Row heights - keep the 19.95pt custom height through row 16
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 16; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.95
}

# ---------------------------------------------------------------------------
# 3. Header row (row 1) values
# ---------------------------------------------------------------------------
$headers = @("exp", "Dataset", "GPU", "GPU_Memory", "GPU_NUM", "batch_size", "Epoch", "memory-Usage", "FLOPs / 1e9", "CPU_US", "CPU_SY", "Load average", "time")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 4. Data rows 2-6
# ---------------------------------------------------------------------------
$data = @(
    @("FR-R50",  "coco2014", "RTX 3070", "8G",  1, 2, 12, $null,              3.8, $null, $null, $null,              "1D16H"),
    @("FR-R101", "coco2014", "RTX 3070", "8G",  1, 2, 12, $null,              7.6, $null, $null, $null,              "2D16H"),
    @("FR-R101", "coco2014", "RTX 3060", "12G", 2, 2, 12, "10000, 8400",      7.6, 12.7,  4.1,   "2.23, 1.72, 0.86", "1D14H"),
    @("FR-R101", "voc0712",  "RTX 3070", "8G",  1, 2, 12, "~7000+",           7.6, $null, $null, $null,              "22H"),
    @("FR-R101", "voc0712",  "RTX 3060", "12G", 2, 2, 12, "7392, 6490",       7.6, 13,    4.2,   "2.09, 1.82, 1.09", "16H")
)

for ($rIdx = 0; $rIdx -lt $data.Length; $rIdx++) {
    $rowNum = $rIdx + 2
    $rowVals = $data[$rIdx]
    for ($cIdx = 0; $cIdx -lt $rowVals.Length; $cIdx++) {
        $val = $rowVals[$cIdx]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $cIdx + 1).Value = $val
        }
    }
}

# ---------------------------------------------------------------------------
# 5. Fonts / alignment
#    header (row1): bold, size 11, 微软雅黑, centered
#    body   (rows2-6): size 12 (not bold), 微软雅黑, centered
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:M1")
$bodyRange   = $ws.Range("A2:M6")

$headerRange.Font.Name = "微软雅黑"
$bodyRange.Font.Name   = "微软雅黑"

$headerRange.Font.Size = 11
$bodyRange.Font.Size   = 12

$headerRange.Font.Bold = $true
$bodyRange.Font.Bold   = $false

$headerRange.HorizontalAlignment = -4108
$bodyRange.HorizontalAlignment   = -4108

$headerRange.VerticalAlignment = -4108
$bodyRange.VerticalAlignment   = -4108

# ---------------------------------------------------------------------------
# 6. Selection
# ---------------------------------------------------------------------------
$ws.Range("G11").Select()
